$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Title paragraph: "Dor Cohen - Full Stack Developer" loses bold (w:b/w:bCs)
#    from the paragraph mark run properties and every run's rPr.
# ---------------------------------------------------------------------------
$xmlTitle = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="009E42E9" w:rsidRPr="00424219" w:rsidRDefault="001F069E" w:rsidP="009E42E9"><w:pPr><w:rPr><w:color w:val="4472C4" w:themeColor="accent5"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00424219"><w:rPr><w:color w:val="4472C4" w:themeColor="accent5"/></w:rPr><w:t>Dor</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00424219"><w:rPr><w:color w:val="4472C4" w:themeColor="accent5"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="002522AF" w:rsidRPr="00424219"><w:rPr><w:color w:val="4472C4" w:themeColor="accent5"/></w:rPr><w:t>Cohen</w:t></w:r><w:r w:rsidRPr="00424219"><w:rPr><w:color w:val="4472C4" w:themeColor="accent5"/></w:rPr><w:t xml:space="preserve"> &#8211; Full Stack Developer</w:t></w:r></w:p>'
$d.Paragraphs(1).Range.InsertXML($xmlTitle)

# ---------------------------------------------------------------------------
# 2) "Phone: 053-7171650" bullet: drop the stray Hyperlink rPr on the pPr.
# ---------------------------------------------------------------------------
$xmlPhone = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00291FE8" w:rsidRPr="00291FE8" w:rsidRDefault="00291FE8" w:rsidP="00291FE8"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr></w:pPr><w:r><w:t>Phone: 053-7171650</w:t></w:r></w:p>'
$d.Paragraphs(5).Range.InsertXML($xmlPhone)

# ---------------------------------------------------------------------------
# 3) Summary paragraph: reworded, with the tech-stack sentence and the
#    _GoBack bookmark relocated here.
# ---------------------------------------------------------------------------
$xmlSummary = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="009E42E9" w:rsidRPr="009E42E9" w:rsidRDefault="009E42E9" w:rsidP="00A97846"><w:r w:rsidRPr="00424219"><w:rPr><w:b/><w:bCs/><w:color w:val="4472C4" w:themeColor="accent5"/></w:rPr><w:t>Summary</w:t></w:r><w:r w:rsidRPr="00424219"><w:rPr><w:color w:val="4472C4" w:themeColor="accent5"/></w:rPr><w:t xml:space="preserve">: </w:t></w:r><w:r w:rsidRPr="009E42E9"><w:t xml:space="preserve">A motivated Full Stack Developer with a strong background in digital marketing. Graduated from Coding Academy&#8217;s Full Stack Development </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="009E42E9"><w:t>Bootcamp</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="009E42E9"><w:t xml:space="preserve">. </w:t></w:r><w:r w:rsidRPr="009E42E9"><w:t xml:space="preserve">Proficient in </w:t></w:r><w:r w:rsidRPr="009E42E9"><w:t xml:space="preserve">tech-stack: React, Node.JS, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="009E42E9"><w:t>mongoDB</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="009E42E9"><w:t>, HTML, CSS, SASS, JavaScript, jQuery, PHP</w:t></w:r><w:r w:rsidRPr="009E42E9"><w:t>, REST API</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r w:rsidRPr="009E42E9"><w:t>.</w:t></w:r><w:r w:rsidRPr="009E42E9"><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidRPr="009E42E9"><w:t>Problem-solving</w:t></w:r><w:r w:rsidRPr="009E42E9"><w:t>, attention to detail, and effective communication and teamwork.</w:t></w:r></w:p>'
$d.Paragraphs(7).Range.InsertXML($xmlSummary)

# ---------------------------------------------------------------------------
# 4) "Digital Marketing & WordPress developer" heading becomes
#    "Digital Marketing Manager & WordPress".
# ---------------------------------------------------------------------------
$rpr = '<w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:u w:val="single"/></w:rPr>'
$xmlHeading = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="009E42E9" w:rsidRPr="00424219" w:rsidRDefault="009E42E9" w:rsidP="00CD2FA3"><w:pPr>' + $rpr + '</w:pPr><w:r w:rsidRPr="00424219">' + $rpr + '<w:t>Digital Marketing</w:t></w:r><w:r w:rsidR="00CD2FA3" w:rsidRPr="00424219">' + $rpr + '<w:t xml:space="preserve"> Manager </w:t></w:r><w:r w:rsidR="00CD2FA3" w:rsidRPr="00424219">' + $rpr + '<w:t>&amp;</w:t></w:r><w:r w:rsidR="00CD2FA3" w:rsidRPr="00424219">' + $rpr + '<w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00C438CB" w:rsidRPr="00424219">' + $rpr + '<w:t>WordPress</w:t></w:r><w:r w:rsidRPr="00424219">' + $rpr + '<w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00424219">' + $rpr + '<w:t>Medton-hedim</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00424219">' + $rpr + '<w:t xml:space="preserve"> (2022-2023):</w:t></w:r></w:p>'
$d.Paragraphs(10).Range.InsertXML($xmlHeading)

# ---------------------------------------------------------------------------
# 5) "Optimized and building websites..." bullet.
# ---------------------------------------------------------------------------
$xmlOptimized = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="009E42E9" w:rsidRPr="009E42E9" w:rsidRDefault="00B131D1" w:rsidP="009E42E9"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="6"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Optimized and building </w:t></w:r><w:r><w:t xml:space="preserve">new </w:t></w:r><w:r w:rsidR="009E42E9" w:rsidRPr="009E42E9"><w:t>websites and landing pages using WordPress.</w:t></w:r></w:p>'
$d.Paragraphs(12).Range.InsertXML($xmlOptimized)

# ---------------------------------------------------------------------------
# 6) "Conducted A/B tests..." bullet.
# ---------------------------------------------------------------------------
$xmlAbTests = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="009E42E9" w:rsidRPr="009E42E9" w:rsidRDefault="009E42E9" w:rsidP="009E42E9"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="6"/></w:numPr></w:pPr><w:r w:rsidRPr="009E42E9"><w:t>Conducted A/B tests to improve performance</w:t></w:r><w:r w:rsidRPr="009E42E9"><w:t xml:space="preserve">, </w:t></w:r><w:r w:rsidRPr="009E42E9"><w:t xml:space="preserve">goals of the company&#8217;s </w:t></w:r><w:r w:rsidRPr="009E42E9"><w:t>KPI</w:t></w:r><w:r w:rsidRPr="009E42E9"><w:t>.</w:t></w:r></w:p>'
$d.Paragraphs(13).Range.InsertXML($xmlAbTests)

# ---------------------------------------------------------------------------
# 7) "Analyzed data using Excel..." bullet.
# ---------------------------------------------------------------------------
$xmlAnalyzed = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="009E42E9" w:rsidRPr="009E42E9" w:rsidRDefault="009E42E9" w:rsidP="009E42E9"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="6"/></w:numPr></w:pPr><w:r w:rsidRPr="009E42E9"><w:t xml:space="preserve">Analyzed data using </w:t></w:r><w:r w:rsidRPr="009E42E9"><w:t>google analytics, google search console, tag manager, excel, screaming frog</w:t></w:r><w:r w:rsidRPr="009E42E9"><w:t>.</w:t></w:r></w:p>'
$d.Paragraphs(14).Range.InsertXML($xmlAnalyzed)

# ---------------------------------------------------------------------------
# 8) "Proficient in on-site and off-site SEO promotion." -> "SEO on-site and off-site"
# ---------------------------------------------------------------------------
$d.Paragraphs(15).Range.Find.Execute("Proficient in on-site and off-site SEO promotion.", $true, $false, $false, $false, $false, $true, 1, $false, "SEO on-site and off-site", 2)

# ---------------------------------------------------------------------------
# 9) "Provided top-notch support to clients..." bullet.
# ---------------------------------------------------------------------------
$xmlProvided = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="009E42E9" w:rsidRPr="009E42E9" w:rsidRDefault="009E42E9" w:rsidP="009E42E9"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr></w:pPr><w:r w:rsidRPr="009E42E9"><w:t xml:space="preserve">Provided </w:t></w:r><w:r w:rsidRPr="009E42E9"><w:t>support to clients, meeting goals effectively.</w:t></w:r></w:p>'
$d.Paragraphs(19).Range.InsertXML($xmlProvided)

# ---------------------------------------------------------------------------
# 10) "Built successful e-commerce websites..." bullet.
# ---------------------------------------------------------------------------
$xmlShopify = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="009E42E9" w:rsidRPr="009E42E9" w:rsidRDefault="009E42E9" w:rsidP="009E42E9"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="8"/></w:numPr></w:pPr><w:r w:rsidRPr="009E42E9"><w:t xml:space="preserve">Built </w:t></w:r><w:r w:rsidRPr="009E42E9"><w:t xml:space="preserve">e-commerce websites using </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="009E42E9"><w:t>shopify</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="009E42E9"><w:t>.</w:t></w:r></w:p>'
$d.Paragraphs(22).Range.InsertXML($xmlShopify)

# ---------------------------------------------------------------------------
# 11) "Digital Marketing Course (6 months), HackerU College(2021)" bullet:
#     drop the stray lastRenderedPageBreak.
# ---------------------------------------------------------------------------
$xmlCourse = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00050956" w:rsidRDefault="00050956" w:rsidP="00050956"><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="10"/></w:numPr></w:pPr><w:r w:rsidRPr="00CD2FA3"><w:t xml:space="preserve">Digital Marketing Course (6 months), </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00CD2FA3"><w:t>HackerU</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00CD2FA3"><w:t xml:space="preserve"> College</w:t></w:r><w:r><w:t>(2021)</w:t></w:r></w:p>'
$d.Paragraphs(27).Range.InsertXML($xmlCourse)

# ---------------------------------------------------------------------------
# 12) Remove the whole trailing Skills / Military Experience / Contact
#     Information block (paragraphs 29..46), keeping the blank paragraph
#     (28) that precedes the section break.
# ---------------------------------------------------------------------------
$startPar = $d.Paragraphs(29)
$endPar = $d.Paragraphs($d.Paragraphs.Count)
$r = $d.Range($startPar.Range.Start, $endPar.Range.End)
$r.Delete()

# ---------------------------------------------------------------------------
# 13) Remove the blank paragraph that used to sit between the Summary and
#     "Professional Experience:" heading.
# ---------------------------------------------------------------------------
$d.Paragraphs(8).Range.Delete()
